$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 14: AN (Art of Negotiation) ---
$ws.Cells.Item(14, 1).Value = "7MP501NE22"
$ws.Cells.Item(14, 2).Value = "AN"
$ws.Cells.Item(14, 3).Value = "Art of Negotiation"
$ws.Cells.Item(14, 4).Value = "Prof. Nitin Pillai"
$ws.Cells.Item(14, 5).Value = "T6"
$ws.Cells.Item(14, 6).Value = "A,B"

# --- New row 15: DC (Digital Consulting), second faculty/section ---
$ws.Cells.Item(15, 1).Value = "7MP708SE22"
$ws.Cells.Item(15, 2).Value = "DC"
$ws.Cells.Item(15, 3).Value = "Digital Consulting"
$ws.Cells.Item(15, 4).Value = "Prof. Sapan Oza (VF) "
$ws.Cells.Item(15, 5).Value = "T6"

# --- Copy the "Venue" column cell format (centered / shaded / bordered / wrapped) onto the new Venue cells ---
$ws.Range("E13").Copy()
$ws.Range("E14").PasteSpecial(-4122)
$ws.Range("E13").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Match the natural (content-driven) row heights used by the rest of the sheet ---
$ws.Rows.Item(14).RowHeight = 15.6
$ws.Rows.Item(15).RowHeight = 15.6

# --- Update view: scroll so row 10 is the top visible row, select the newly added Venue cell ---
$ws.Range("E15").Select()

Write-Host "done"
